# Generate Report for handback
# Two source files ("70ddc363-..." and "7593a10e-...") have come back from
# the translators in sync with en-US, so:
#  - their Status flips from "Ready for handoff" to
#    "Handed back: in sync with en-US" (Overview + both language sheets)
#  - the "Latest Target File" / "Latest Handback File" columns get filled in
#    (mirroring the source .md / handoff .xlf, since they are in sync)
#  - "Latest Handback DateTime" gets a real timestamp instead of the epoch
#    placeholder

$wb = $excel.ActiveWorkbook

$newStatus = "Handed back: in sync with en-US"

# ---------------------------------------------------------------------
# Overview sheet: just the status text changes (shared by both languages)
# ---------------------------------------------------------------------
$overview = $wb.Worksheets.Item("Overview")
$overview.Range("B2").Value = $newStatus
$overview.Range("C2").Value = $newStatus
$overview.Range("B3").Value = $newStatus
$overview.Range("C3").Value = $newStatus

# ---------------------------------------------------------------------
# zh-cn sheet
# ---------------------------------------------------------------------
$zh = $wb.Worksheets.Item("zh-cn")

$zh.Range("B2").Value = $newStatus
$zh.Range("B3").Value = $newStatus

$zh.Range("E2").Value = "70ddc363-8daf-465a-8201-3d31f936189c.md"
$zh.Hyperlinks.Add($zh.Range("E2"), "https://github.com/OpenLocalizationTest/oltest/blob/5586e1488afb2ebe7bc27707f6a08e46136f998d/e2e/70ddc363-8daf-465a-8201-3d31f936189c.md", "", "", "70ddc363-8daf-465a-8201-3d31f936189c.md") | Out-Null

$zh.Range("F2").Value = "70ddc363-8daf-465a-8201-3d31f936189c.d9077f3a288e5ecf56dbde5728061a4f845f8330.zh-cn.xlf"
$zh.Hyperlinks.Add($zh.Range("F2"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/32b47fff5b5d8aa934277faeb8876fd29304158d/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/yuwzho/70ddc363-8daf-465a-8201-3d31f936189c.d9077f3a288e5ecf56dbde5728061a4f845f8330.zh-cn.xlf", "", "", "70ddc363-8daf-465a-8201-3d31f936189c.d9077f3a288e5ecf56dbde5728061a4f845f8330.zh-cn.xlf") | Out-Null

$zh.Range("E3").Value = "7593a10e-691d-4a59-b935-bd2d1ef4e50c.md"
$zh.Hyperlinks.Add($zh.Range("E3"), "https://github.com/OpenLocalizationTest/oltest/blob/5586e1488afb2ebe7bc27707f6a08e46136f998d/e2e/7593a10e-691d-4a59-b935-bd2d1ef4e50c.md", "", "", "7593a10e-691d-4a59-b935-bd2d1ef4e50c.md") | Out-Null

$zh.Range("F3").Value = "7593a10e-691d-4a59-b935-bd2d1ef4e50c.34df5b93068e463ea8bca512aedee91af2703221.zh-cn.xlf"
$zh.Hyperlinks.Add($zh.Range("F3"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/32b47fff5b5d8aa934277faeb8876fd29304158d/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/yuwzho/7593a10e-691d-4a59-b935-bd2d1ef4e50c.34df5b93068e463ea8bca512aedee91af2703221.zh-cn.xlf", "", "", "7593a10e-691d-4a59-b935-bd2d1ef4e50c.34df5b93068e463ea8bca512aedee91af2703221.zh-cn.xlf") | Out-Null

$zh.Range("G2").Value = "2016-01-13 15:47:43"
$zh.Range("G3").Value = "2016-01-13 15:47:43"

# ---------------------------------------------------------------------
# de-de sheet
# ---------------------------------------------------------------------
$de = $wb.Worksheets.Item("de-de")

$de.Range("B2").Value = $newStatus
$de.Range("B3").Value = $newStatus

$de.Range("E2").Value = "70ddc363-8daf-465a-8201-3d31f936189c.md"
$de.Hyperlinks.Add($de.Range("E2"), "https://github.com/OpenLocalizationTest/oltest/blob/5586e1488afb2ebe7bc27707f6a08e46136f998d/e2e/70ddc363-8daf-465a-8201-3d31f936189c.md", "", "", "70ddc363-8daf-465a-8201-3d31f936189c.md") | Out-Null

$de.Range("F2").Value = "70ddc363-8daf-465a-8201-3d31f936189c.d9077f3a288e5ecf56dbde5728061a4f845f8330.de-de.xlf"
$de.Hyperlinks.Add($de.Range("F2"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/3bb37d8cc09e0b4a18c2f48c7b94ea7796593964/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/yuwzho/70ddc363-8daf-465a-8201-3d31f936189c.d9077f3a288e5ecf56dbde5728061a4f845f8330.de-de.xlf", "", "", "70ddc363-8daf-465a-8201-3d31f936189c.d9077f3a288e5ecf56dbde5728061a4f845f8330.de-de.xlf") | Out-Null

$de.Range("E3").Value = "7593a10e-691d-4a59-b935-bd2d1ef4e50c.md"
$de.Hyperlinks.Add($de.Range("E3"), "https://github.com/OpenLocalizationTest/oltest/blob/5586e1488afb2ebe7bc27707f6a08e46136f998d/e2e/7593a10e-691d-4a59-b935-bd2d1ef4e50c.md", "", "", "7593a10e-691d-4a59-b935-bd2d1ef4e50c.md") | Out-Null

$de.Range("F3").Value = "7593a10e-691d-4a59-b935-bd2d1ef4e50c.34df5b93068e463ea8bca512aedee91af2703221.de-de.xlf"
$de.Hyperlinks.Add($de.Range("F3"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/3bb37d8cc09e0b4a18c2f48c7b94ea7796593964/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/yuwzho/7593a10e-691d-4a59-b935-bd2d1ef4e50c.34df5b93068e463ea8bca512aedee91af2703221.de-de.xlf", "", "", "7593a10e-691d-4a59-b935-bd2d1ef4e50c.34df5b93068e463ea8bca512aedee91af2703221.de-de.xlf") | Out-Null

$de.Range("G2").Value = "2016-01-13 15:48:03"
$de.Range("G3").Value = "2016-01-13 15:48:03"

Write-Output "done"
